$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'304.72"
$ws.Range("E2").Value = "'-4.92%"
$ws.Range("E3").Value = "'-7.17%"
$ws.Range("D4").Value = "'5.041"
$ws.Range("E4").Value = "'-2.81%"
$ws.Range("D5").Value = "'0.07669"
$ws.Range("E5").Value = "'-5.81%"
$ws.Range("D6").Value = "'4.250"
$ws.Range("E6").Value = "'-2.14%"
$ws.Range("D7").Value = "'1.586"
$ws.Range("E7").Value = "'-11.39%"
$ws.Range("D8").Value = "'0.8819"
$ws.Range("E8").Value = "'-7.21%"
$ws.Range("D9").Value = "'0.09910"
$ws.Range("E9").Value = "'-11.18%"
$ws.Range("D10").Value = "'0.1722"
$ws.Range("E10").Value = "'-6.65%"
$ws.Range("D11").Value = "'0.08950"
$ws.Range("E11").Value = "'-4.28%"
$ws.Range("D12").Value = "'0.04434"
$ws.Range("E12").Value = "'-5.02%"
$ws.Range("D13").Value = "'0.1054"
$ws.Range("E13").Value = "'-0.64%"
$ws.Range("D14").Value = "'0.001282"
$ws.Range("E14").Value = "'-0.54%"
$ws.Range("D15").Value = "'0.005829"
$ws.Range("E15").Value = "'-1.46%"
$ws.Range("D16").Value = "'3.360"
$ws.Range("E16").Value = "'-0.16%"
$ws.Range("D17").Value = "'2.419"
$ws.Range("E17").Value = "'-4.37%"
$ws.Range("E18").Value = "'-0.10%"
$ws.Range("D19").Value = "'7.068"
$ws.Range("E19").Value = "'-5.08%"
$ws.Range("D20").Value = "'0.1351"
$ws.Range("E20").Value = "'-3.06%"
$ws.Range("D21").Value = "'0.3238"
$ws.Range("E21").Value = "'23.22%"
$ws.Range("D22").Value = "'0.04206"
$ws.Range("E22").Value = "'0.42%"
$ws.Range("D23").Value = "'0.001196"
$ws.Range("E23").Value = "'-4.56%"
$ws.Range("D24").Value = "'0.004057"
$ws.Range("E24").Value = "'-6.21%"
$ws.Range("D25").Value = "'0.0001225"
$ws.Range("E25").Value = "'10.23%"
$ws.Range("E26").Value = "'-0.34%"
$ws.Range("D38").Value = "'0.02335"
$ws.Range("E38").Value = "'-9.65%"
$ws.Range("D39").Value = "'0.05118"
$ws.Range("E39").Value = "'-7.34%"
$ws.Range("D40").Value = "'0.007952"
$ws.Range("E40").Value = "'3.30%"
$ws.Range("E41").Value = "'-5.17%"
$ws.Range("D42").Value = "'0.006656"
$ws.Range("E42").Value = "'0.54%"
$ws.Range("D43").Value = "'0.001992"
$ws.Range("E43").Value = "'-5.88%"
$ws.Range("D44").Value = "'0.008502"
$ws.Range("E44").Value = "'0.69%"
$ws.Range("D45").Value = "'0.3021"
$ws.Range("E45").Value = "'-12.92%"
$ws.Range("D46").Value = "'0.00006539"
$ws.Range("E46").Value = "'-6.47%"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'0.07%"
$ws.Range("D48").Value = "'0.007013"
$ws.Range("E48").Value = "'97.82%"
$ws.Range("D49").Value = "'0.003384"
$ws.Range("E49").Value = "'-2.83%"
$ws.Range("D50").Value = "'0.00002104"
$ws.Range("E50").Value = "'0.07%"
$ws.Range("D51").Value = "'0.0002004"
$ws.Range("E51").Value = "'0.07%"
